$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "42.239.45"
$ws.Range("E2").Value = "  -1.06%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.244.66"
$ws.Range("E3").Value = "  -1.03%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - BNB
$ws.Range("D5").Value = "247.04"
$ws.Range("E5").Value = "  -0.76%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.621"
$ws.Range("E6").Value = "  -2.42%  "

# Row 7
$ws.Range("D7").Value = "74.88"
$ws.Range("E7").Value = "  -3.00%  "

# Row 8
$ws.Range("E8").Value = "  +0.08%  "

# Row 9
$ws.Range("D9").Value = "0.619"
$ws.Range("E9").Value = "  -3.77%  "

# Row 10
$ws.Range("D10").Value = "42.25"
$ws.Range("E10").Value = "  +6.17%  "

# Row 11
$ws.Range("D11").Value = "0.0944"
$ws.Range("E11").Value = "  -2.13%  "

# Row 12
$ws.Range("D12").Value = "7.16"
$ws.Range("E12").Value = "  -0.79%  "

# Row 13
$ws.Range("D13").Value = "0.103"
$ws.Range("E13").Value = "  -2.41%  "

# Row 14
$ws.Range("D14").Value = "14.53"
$ws.Range("E14").Value = "  -2.93%  "

# Row 15
$ws.Range("D15").Value = "0.854"
$ws.Range("E15").Value = "  -0.92%  "

# Row 16
$ws.Range("D16").Value = "2.256.66"
$ws.Range("E16").Value = "  -0.63%  "

# Row 17
$ws.Range("D17").Value = "42.114.01"
$ws.Range("E17").Value = "  -1.15%  "

# Row 18
$ws.Range("D18").Value = "0.0₃0987"
$ws.Range("E18").Value = "  +0.07%  "

# Row 19
$ws.Range("E19").Value = "  -0.04%  "

# Row 20
$ws.Range("D20").Value = "72.23"
$ws.Range("E20").Value = "  +0.54%  "

# Row 21 - now ImmutableX (was BitcoinCash)
$ws.Range("B21").Value = "ImmutableX"
$ws.Range("C21").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D21").Value = "2.22"
$ws.Range("E21").Value = "  +4.93%  "

# Row 22 - now BitcoinCash (was ImmutableX)
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "231.95"
$ws.Range("E22").Value = "  -1.19%  "

# Row 23
$ws.Range("D23").Value = "9.04"
$ws.Range("E23").Value = "  +42.02%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.30%  "

# Row 25
$ws.Range("D25").Value = "11.38"
$ws.Range("E25").Value = "  +0.91%  "

# Row 26
$ws.Range("E26").Value = "  -4.60%  "

# Row 27
$ws.Range("E27").Value = "  -2.20%  "

# Row 28
$ws.Range("D28").Value = "2.23"
$ws.Range("E28").Value = "  +2.69%  "

# Row 29
$ws.Range("D29").Value = "169.41"
$ws.Range("E29").Value = "  +1.19%  "

# Row 30
$ws.Range("D30").Value = "20.69"
$ws.Range("E30").Value = "  -0.75%  "

# Row 31
$ws.Range("E31").Value = "  -3.70%  "

# Row 32
$ws.Range("D32").Value = "31.04"
$ws.Range("E32").Value = "  +1.49%  "

# Row 33
$ws.Range("E33").Value = "  -1.56%  "

# Row 34
$ws.Range("E34").Value = "  -1.14%  "

# Row 35
$ws.Range("D35").Value = "5.21"
$ws.Range("E35").Value = "  +11.20%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.50"
$ws.Range("E36").Value = "  -0.93%  "

# Row 37
$ws.Range("D37").Value = "0.0312"
$ws.Range("E37").Value = "  +2.61%  "

# Row 38
$ws.Range("D38").Value = "13.79"
$ws.Range("E38").Value = "  +0.31%  "

# Row 39
$ws.Range("E39").Value = "  -2.73%  "

# Row 40
$ws.Range("D40").Value = "5.76"
$ws.Range("E40").Value = "  -0.90%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "62.30"
$ws.Range("E41").Value = "  +2.43%  "

# Row 42
$ws.Range("D42").Value = "0.204"
$ws.Range("E42").Value = "  -1.26%  "

# Row 43
$ws.Range("D43").Value = "106.31"
$ws.Range("E43").Value = "  -3.22%  "

# Row 44
$ws.Range("D44").Value = "0.103"
$ws.Range("E44").Value = "  +2.56%  "

# Row 45
$ws.Range("E45").Value = "  -1.72%  "

# Row 46
$ws.Range("D46").Value = "0.996"
$ws.Range("E46").Value = "  -0.34%  "

# Row 47
$ws.Range("E47").Value = "  -2.21%  "

# Row 48
$ws.Range("E48").Value = "  +0.65%  "

# Row 49
$ws.Range("E49").Value = "  -5.88%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.30"
$ws.Range("E50").Value = "  +3.29%  "

# Row 51 - now SynthetixNetwork (was BitTorrent-New)
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").Value = "4.12"
$ws.Range("E51").Value = "  -2.16%  "
